# Updates the crypto price/volume/hour table to the latest GitHub Actions
# scrape (Mon Dec 19 21:04:43 UTC 2022). Rows 10-18 also get re-ranked
# (coins move up/down the Best/Worst-in-24h board), and a few "Bestin24h"
# / "Worstin24h" badges move to their new rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and G (Hora) hold numeric-looking values that are
# stored as literal text in the workbook (e.g. "0.0001499" must keep its
# trailing digits, not become the number 0.0001499 with different
# formatting). Mark the cells we are about to touch as Text first so
# Excel does not silently convert the assigned strings to numbers.
$ws.Range("D2:D20").NumberFormat = "@"
$ws.Range("D22:D24").NumberFormat = "@"
$ws.Range("D40:D45").NumberFormat = "@"
$ws.Range("D47:D50").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Write the new cell values.
# Row 2
$ws.Range("D2").Value = "243.49"
$ws.Range("G2").Value = "21"
# Row 3
$ws.Range("D3").Value = "21.49"
$ws.Range("G3").Value = "21"
# Row 4
$ws.Range("D4").Value = "5.224"
$ws.Range("G4").Value = "21"
# Row 5
$ws.Range("D5").Value = "0.05603"
$ws.Range("G5").Value = "21"
# Row 6
$ws.Range("D6").Value = "3.364"
$ws.Range("G6").Value = "21"
# Row 7
$ws.Range("D7").Value = "6.376"
$ws.Range("G7").Value = "21"
# Row 8
$ws.Range("D8").Value = "0.8049"
$ws.Range("G8").Value = "21"
# Row 9
$ws.Range("D9").Value = "0.9486"
$ws.Range("G9").Value = "21"
# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1439"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("G10").Value = "21"
# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07286"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("G11").Value = "21"
# Row 12
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03151"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G12").Value = "21"
# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03092"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("G13").Value = "21"
# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09277"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("G14").Value = "21"
# Row 15
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "3.565"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("G15").Value = "21"
# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001651"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("G16").Value = "21"
# Row 17
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04702"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("G17").Value = "21"
# Row 18
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005745"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("G18").Value = "21"
# Row 19
$ws.Range("D19").Value = "0.006350"
$ws.Range("G19").Value = "21"
# Row 20
$ws.Range("D20").Value = "0.004991"
$ws.Range("G20").Value = "21"
# Row 21
$ws.Range("G21").Value = "21"
# Row 22
$ws.Range("D22").Value = "0.0001499"
$ws.Range("G22").Value = "21"
# Row 23
$ws.Range("D23").Value = "0.0003097"
$ws.Range("G23").Value = "21"
# Row 24
$ws.Range("D24").Value = "3.756"
$ws.Range("G24").Value = "21"
# Row 25
$ws.Range("G25").Value = "21"
# Row 26
$ws.Range("G26").Value = "21"
# Row 27
$ws.Range("G27").Value = "21"
# Row 28
$ws.Range("G28").Value = "21"
# Row 29
$ws.Range("G29").Value = "21"
# Row 30
$ws.Range("G30").Value = "21"
# Row 31
$ws.Range("G31").Value = "21"
# Row 32
$ws.Range("G32").Value = "21"
# Row 33
$ws.Range("G33").Value = "21"
# Row 34
$ws.Range("G34").Value = "21"
# Row 35
$ws.Range("G35").Value = "21"
# Row 36
$ws.Range("G36").Value = "21"
# Row 37
$ws.Range("G37").Value = "21"
# Row 38
$ws.Range("G38").Value = "21"
# Row 39
$ws.Range("G39").Value = "21"
# Row 40
$ws.Range("D40").Value = "0.03921"
$ws.Range("G40").Value = "21"
# Row 41
$ws.Range("D41").Value = "0.006889"
$ws.Range("G41").Value = "21"
# Row 42
$ws.Range("D42").Value = "0.003399"
$ws.Range("G42").Value = "21"
# Row 43
$ws.Range("D43").Value = "0.1034"
$ws.Range("G43").Value = "21"
# Row 44
$ws.Range("D44").Value = "0.007522"
$ws.Range("G44").Value = "21"
# Row 45
$ws.Range("D45").Value = "0.00005926"
$ws.Range("G45").Value = "21"
# Row 46
$ws.Range("G46").Value = "21"
# Row 47
$ws.Range("D47").Value = "0.0005495"
$ws.Range("E47").Value = "46ACDXExchangeACXT"
$ws.Range("G47").Value = "21"
# Row 48
$ws.Range("D48").Value = "0.6818"
$ws.Range("G48").Value = "21"
# Row 49
$ws.Range("D49").Value = "0.07929"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"
$ws.Range("G49").Value = "21"
# Row 50
$ws.Range("D50").Value = "0.00002099"
$ws.Range("G50").Value = "21"
# Row 51
$ws.Range("G51").Value = "21"
